$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.004.46'
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('D3').Value = '2.914.02'
$ws.Range('E3').Value = '  +0.11%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '''588.92'
$ws.Range('E5').Value = '  +0.29%  '
$ws.Range('D6').Value = '''144.86'
$ws.Range('E6').Value = '  +0.47%  '
$ws.Range('E8').Value = '  +0.22%  '
$ws.Range('D9').Value = '''6.91'
$ws.Range('E9').Value = '  +3.98%  '
$ws.Range('E11').Value = '  -1.51%  '
$ws.Range('E12').Value = '  -0.22%  '
$ws.Range('D13').Value = '''33.42'
$ws.Range('E13').Value = '  +0.16%  '
$ws.Range('E14').Value = '  -0.09%  '
$ws.Range('D15').Value = '3.396.42'
$ws.Range('E15').Value = '  +0.03%  '
$ws.Range('D16').Value = '60.907.50'
$ws.Range('E16').Value = '  +0.23%  '
$ws.Range('E17').Value = '  -0.34%  '
$ws.Range('D18').Value = '2.912.82'
$ws.Range('E18').Value = '  +0.04%  '
$ws.Range('D19').Value = '''432.72'
$ws.Range('E19').Value = '  +1.22%  '
$ws.Range('E20').Value = '  -1.28%  '
$ws.Range('D21').Value = '0.675'
$ws.Range('E21').Value = '  -0.67%  '
$ws.Range('D22').Value = '''7.13'
$ws.Range('E22').Value = '  +0.57%  '
$ws.Range('D23').Value = '''81.42'
$ws.Range('E23').Value = '  +1.24%  '
$ws.Range('E24').Value = '  +0.28%  '
$ws.Range('E25').Value = '  -1.53%  '
$ws.Range('D26').Value = '11.77'
$ws.Range('E26').Value = '  -1.15%  '
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('D28').Value = '2.27'
$ws.Range('E28').Value = '  +5.14%  '
$ws.Range('E29').Value = '  -0.62%  '
$ws.Range('D30').Value = '6.96'
$ws.Range('E30').Value = '  -3.42%  '
$ws.Range('D31').Value = '26.48'
$ws.Range('E31').Value = '  +0.22%  '
$ws.Range('D32').Value = '0.108'
$ws.Range('E32').Value = '  +2.69%  '
$ws.Range('E33').Value = '  -0.09%  '
$ws.Range('D34').Value = '0.0₃0870'
$ws.Range('E34').Value = '  -0.38%  '
$ws.Range('E35').Value = '  +0.17%  '
$ws.Range('E36').Value = '  +0.29%  '
$ws.Range('D37').Value = '''3.00'
$ws.Range('E37').Value = '  -0.27%  '
$ws.Range('D38').Value = '1.98'
$ws.Range('E38').Value = '  -0.78%  '
$ws.Range('E39').Value = '  -3.95%  '
$ws.Range('E40').Value = '  -0.36%  '
$ws.Range('D41').Value = '''0.284'
$ws.Range('E41').Value = '  -4.31%  '
$ws.Range('D42').Value = '''40.91'
$ws.Range('E42').Value = '  -1.15%  '
$ws.Range('D43').Value = '''376.41'
$ws.Range('E43').Value = '  -0.46%  '
$ws.Range('D44').Value = '''0.0345'
$ws.Range('E44').Value = '  -1.29%  '
$ws.Range('D45').Value = '2.692.47'
$ws.Range('E45').Value = '  +0.18%  '
$ws.Range('D46').Value = '''133.64'
$ws.Range('E47').Value = '  -0.04%  '
$ws.Range('D48').Value = '23.75'
$ws.Range('E48').Value = '  -2.47%  '
$ws.Range('E49').Value = '  -0.49%  '
$ws.Range('E50').Value = '  -2.12%  '
$ws.Range('E51').Value = '  -0.39%  '
